$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.186.22'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '1.892.11'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.82'
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5220'
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3750'
$ws.Range("E8").Value = '  -1.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07254'
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.16'
$ws.Range("E10").Value = '  -0.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8970'
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08154'
$ws.Range("E12").Value = '  +6.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '96.62'
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("D14").Value = '1.894.77'
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.265'
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008576'
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.52'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = '27.227.06'
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.074'
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.67'
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.394'
$ws.Range("E23").Value = '  -0.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '147.40'
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.285'
$ws.Range("E25").Value = '  -1.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.743'
$ws.Range("E26").Value = '  +0.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.17'
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '114.84'
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.787'
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09218'
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05031'
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7901'
$ws.Range("E33").Value = '  +0.70%  '
$ws.Range("E34").Value = '  -2.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.431'
$ws.Range("E35").Value = '  +3.76%  '
$ws.Range("E36").Value = '  -1.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.561'
$ws.Range("E37").Value = '  -1.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5638'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01982'
$ws.Range("E39").Value = '  -0.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.073'
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.914'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.523'
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '115.13'
$ws.Range("E43").Value = '  -2.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1513'
$ws.Range("E44").Value = '  -0.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4855'
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.03'
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.614'
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.05'
$ws.Range("E49").Value = '  +1.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.24'
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05937'
$ws.Range("E51").Value = '  +0.14%  '
